# Update "historico" / Palmeiras rodada_26 sheet: replace game-id (A) and
# round/matchweek (E) values to reflect the upcoming rodada 27 poisson_naive run.
# Column E switches from a text label ("Matchweek N") to a plain number N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (A value, E value)
$updates = @{
    2  = @(1, 1)
    3  = @(4, 3)
    4  = @(7, 5)
    5  = @(8, 6)
    6  = @(12, 9)
    7  = @(16, 12)
    8  = @(19, 14)
    9  = @(21, 16)
    10 = @(26, 19)
    11 = @(29, 21)
    12 = @(32, 23)
    13 = @(37, 26)
    14 = @(32, 25)
    15 = @(23, 24)
    16 = @(24, 18)
    17 = @(18, 13)
    18 = @(15, 8)
    19 = @(14, 10)
    20 = @(19, 20)
    21 = @(20, 15)
    22 = @(32, 22)
    23 = @(9, 7)
    24 = @(10, 11)
    25 = @(1, 2)
    26 = @(6, 4)
    27 = @(23, 17)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]   # column A
    $ws.Cells.Item($row, 5).Value = $vals[1]   # column E
}
